$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header column F1 ("REX_DEF"), matching the formatting of the
# existing header cells (B1:E1) which use style index 1 (bold, bordered,
# centered/top-aligned).
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "REX_DEF"

$excel.CutCopyMode = $false
